$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet: "Datos" -> "animales_masiva"
$ws.Name = "animales_masiva"

# Remove the obsolete columns (delete right-to-left so positions stay valid):
#   W -> "composicion_racial"
#   V -> "numero_hierros"
#   M -> "sector"
$ws.Range("W1").EntireColumn.Delete()
$ws.Range("V1").EntireColumn.Delete()
$ws.Range("M1").EntireColumn.Delete()

# Strip the bold/white-on-blue header styling back to the default style
$ws.Range("A1:Z1").ClearFormats()

Write-Host "Done"
